$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3757.5
$ws.Range("J69").Value = 2865.7144
$ws.Range("L69").Value = 8597.143199999999
$ws.Range("N69").Value = -10345.1432
$ws.Range("H70").Value = 1693.9
$ws.Range("I70").Value = 1644.6111
$ws.Range("J70").Value = 1734.2273
$ws.Range("K70").Value = 4933.8333
$ws.Range("L70").Value = 5202.6819
$ws.Range("M70").Value = -4663.8333
$ws.Range("N70").Value = -5742.6819
$ws.Range("H72").Value = 3757.5
$ws.Range("J72").Value = 2865.7144
$ws.Range("L72").Value = 25791.4296
$ws.Range("N72").Value = -34527.4296
$ws.Range("H73").Value = 1693.9
$ws.Range("I73").Value = 1644.6111
$ws.Range("J73").Value = 1734.2273
$ws.Range("K73").Value = 4933.8333
$ws.Range("L73").Value = 5202.6819
$ws.Range("M73").Value = -3997.8333
$ws.Range("N73").Value = -7074.6819
$ws.Range("H76").Value = 3933.1667
$ws.Range("I76").Value = 3554.9656
$ws.Range("K76").Value = 3554.9656
$ws.Range("M76").Value = -3239.9656
$ws.Range("H79").Value = 3933.1667
$ws.Range("I79").Value = 3554.9656
$ws.Range("K79").Value = 3554.9656
$ws.Range("M79").Value = -2462.9656
$ws.Range("H86").Value = 47834.816
$ws.Range("I86").Value = 64691.625
$ws.Range("J86").Value = 2883.3333
$ws.Range("K86").Value = 64691.625
$ws.Range("L86").Value = 2883.3333
$ws.Range("M86").Value = -63568.625
$ws.Range("N86").Value = -5129.3333
$ws.Range("H89").Value = 47834.816
$ws.Range("I89").Value = 64691.625
$ws.Range("J89").Value = 2883.3333
$ws.Range("K89").Value = 323458.125
$ws.Range("L89").Value = 14416.6665
$ws.Range("M89").Value = -317842.125
$ws.Range("N89").Value = -25648.6665
$ws.Range("H137").Value = 880479.75
$ws.Range("I137").Value = 5877.778
$ws.Range("J137").Value = 1151908
$ws.Range("K137").Value = 17633.334
$ws.Range("L137").Value = 3455724
$ws.Range("M137").Value = -15083.334
$ws.Range("N137").Value = -3460824
$ws.Range("H138").Value = 3361.6807
$ws.Range("I138").Value = 1777.5883
$ws.Range("J138").Value = 3851.309
$ws.Range("K138").Value = 5332.7649
$ws.Range("L138").Value = 11553.927
$ws.Range("M138").Value = -192.7649000000001
$ws.Range("N138").Value = -21833.927

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18832.787
$ws.Range("I32").Value = 19840.59
$ws.Range("K32").Value = 19840.59
$ws.Range("M32").Value = -19553.59
$ws.Range("H61").Value = 11367.92
$ws.Range("I61").Value = 8289
$ws.Range("J61").Value = 14703.417
$ws.Range("K61").Value = 8289
$ws.Range("L61").Value = 14703.417
$ws.Range("M61").Value = -8077
$ws.Range("N61").Value = -15127.417
$ws.Range("H74").Value = 3571.524
$ws.Range("I74").Value = 1171.5555
$ws.Range("J74").Value = 17971.334
$ws.Range("K74").Value = 1171.5555
$ws.Range("L74").Value = 17971.334
$ws.Range("M74").Value = -297.5554999999999
$ws.Range("N74").Value = -19719.334
$ws.Range("H77").Value = 3571.524
$ws.Range("I77").Value = 1171.5555
$ws.Range("J77").Value = 17971.334
$ws.Range("K77").Value = 5857.7775
$ws.Range("L77").Value = 89856.67
$ws.Range("M77").Value = -1489.7775
$ws.Range("N77").Value = -98592.67
$ws.Range("H102").Value = 3292
$ws.Range("I102").Value = 3167.8572
$ws.Range("J102").Value = 3639.6
$ws.Range("K102").Value = 3167.8572
$ws.Range("L102").Value = 3639.6
$ws.Range("M102").Value = -1545.8572
$ws.Range("N102").Value = -6883.6
$ws.Range("H136").Value = 11367.92
$ws.Range("I136").Value = 8289
$ws.Range("J136").Value = 14703.417
$ws.Range("K136").Value = 24867
$ws.Range("L136").Value = 44110.251
$ws.Range("M136").Value = -22317
$ws.Range("N136").Value = -49210.251

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2868.647
$ws.Range("I99").Value = 897.25
$ws.Range("K99").Value = 897.25
$ws.Range("M99").Value = 600.75
$ws.Range("H134").Value = 52162.85
$ws.Range("I134").Value = 2319.0557
$ws.Range("K134").Value = 6957.1671
$ws.Range("M134").Value = -4422.1671

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 611867.75
$ws.Range("I31").Value = 7179.1333
$ws.Range("J31").Value = 1197050.2
$ws.Range("K31").Value = 7179.1333
$ws.Range("L31").Value = 1197050.2
$ws.Range("M31").Value = -6884.1333
$ws.Range("N31").Value = -1197640.2
$ws.Range("H34").Value = 611867.75
$ws.Range("I34").Value = 7179.1333
$ws.Range("J34").Value = 1197050.2
$ws.Range("K34").Value = 7179.1333
$ws.Range("L34").Value = 1197050.2
$ws.Range("M34").Value = -6977.1333
$ws.Range("N34").Value = -1197454.2
$ws.Range("H58").Value = 2021758.6
$ws.Range("I58").Value = 4133613
$ws.Range("J58").Value = 1724.0869
$ws.Range("K58").Value = 4133613
$ws.Range("L58").Value = 1724.0869
$ws.Range("M58").Value = -4133410
$ws.Range("N58").Value = -2130.0869
$ws.Range("H94").Value = 352.54544
$ws.Range("J94").Value = 367.8
$ws.Range("L94").Value = 367.8
$ws.Range("N94").Value = -1269.8
$ws.Range("H134").Value = 2353.6487
$ws.Range("I134").Value = 2025.0476
$ws.Range("K134").Value = 6075.142800000001
$ws.Range("M134").Value = -3540.142800000001
$ws.Range("H136").Value = 2021758.6
$ws.Range("I136").Value = 4133613
$ws.Range("J136").Value = 1724.0869
$ws.Range("K136").Value = 12400839
$ws.Range("L136").Value = 5172.2607
$ws.Range("M136").Value = -12398289
$ws.Range("N136").Value = -10272.2607

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 988.25
$ws.Range("I122").Value = 730
$ws.Range("J122").Value = 1044.3914
$ws.Range("K122").Value = 6570
$ws.Range("L122").Value = 9399.5226
$ws.Range("M122").Value = -4120
$ws.Range("N122").Value = -14299.5226
$ws.Range("H131").Value = 1183.3036
$ws.Range("J131").Value = 1043.3695
$ws.Range("L131").Value = 3130.1085
$ws.Range("N131").Value = -13210.1085
$ws.Range("H140").Value = 1928.0541
$ws.Range("I140").Value = 1305.3572
$ws.Range("J140").Value = 3865.3333
$ws.Range("K140").Value = 3916.0716
$ws.Range("L140").Value = 11595.9999
$ws.Range("M140").Value = 1263.9284
$ws.Range("N140").Value = -21955.9999
$ws.Range("H141").Value = 1821.3334
$ws.Range("I141").Value = 1821.3334
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5464.0002
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -284.0002000000004
$ws.Range("N141").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 15556.667
$ws.Range("I132").Value = 11085.091
$ws.Range("J132").Value = 27853.5
$ws.Range("K132").Value = 33255.273
$ws.Range("L132").Value = 83560.5
$ws.Range("M132").Value = -30725.273
$ws.Range("N132").Value = -88620.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H16").Value = 6981.125
$ws.Range("I16").Value = 835.5714
$ws.Range("K16").Value = 835.5714
$ws.Range("M16").Value = -665.5714
$ws.Range("H40").Value = 3316.8696
$ws.Range("I40").Value = 3246.9333
$ws.Range("J40").Value = 3448
$ws.Range("K40").Value = 3246.9333
$ws.Range("L40").Value = 3448
$ws.Range("M40").Value = -3110.9333
$ws.Range("N40").Value = -3720
$ws.Range("H82").Value = 2252.25
$ws.Range("I82").Value = 1501.5
$ws.Range("J82").Value = 3003
$ws.Range("K82").Value = 1501.5
$ws.Range("L82").Value = 3003
$ws.Range("M82").Value = -1140.5
$ws.Range("N82").Value = -3725
$ws.Range("H85").Value = 2252.25
$ws.Range("I85").Value = 1501.5
$ws.Range("J85").Value = 3003
$ws.Range("K85").Value = 1501.5
$ws.Range("L85").Value = 3003
$ws.Range("M85").Value = -253.5
$ws.Range("N85").Value = -5499
$ws.Range("H99").Value = 61000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 61000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 61000
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -66990
$ws.Range("H100").Value = 5866.6
$ws.Range("I100").Value = 5857.143
$ws.Range("J100").Value = 5874.875
$ws.Range("K100").Value = 5857.143
$ws.Range("L100").Value = 5874.875
$ws.Range("M100").Value = -5316.143
$ws.Range("N100").Value = -6956.875
$ws.Range("H134").Value = 49812.332
$ws.Range("J134").Value = 55774.8
$ws.Range("L134").Value = 55774.8
$ws.Range("N134").Value = -65914.8
$ws.Range("H140").Value = 78632
$ws.Range("J140").Value = 78632
$ws.Range("L140").Value = 78632
$ws.Range("N140").Value = -88992
$ws.Range("H141").Value = 70000
$ws.Range("J141").Value = 70000
$ws.Range("L141").Value = 70000
$ws.Range("N141").Value = -80360

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1937.5
$ws.Range("I122").Value = 2225
$ws.Range("J122").Value = 1650
$ws.Range("K122").Value = 6675
$ws.Range("L122").Value = 4950
$ws.Range("M122").Value = -4225
$ws.Range("N122").Value = -9850
